$d = $word.ActiveDocument
$d.Content.Find.Execute("Mark, Morrison", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Mark Morrison", 2)
